$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.5
$ws.Range("J3").Value = 3.4
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 2.25
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 1.13
$ws.Range("AF3").Value = 23
$ws.Range("AO3").Value = 15
$ws.Range("O5").Value = 1.7
$ws.Range("AA5").Value = 2.85
$ws.Range("AB5").Value = 1.37
$ws.Range("G9").Value = 1.36
$ws.Range("H9").Value = 4.35
$ws.Range("I9").Value = 8.75
$ws.Range("J9").Value = 1.87
$ws.Range("K9").Value = 2.3
$ws.Range("L9").Value = 7.3
$ws.Range("N9").Value = 8
$ws.Range("O9").Value = 1.23
$ws.Range("P9").Value = 3.75
$ws.Range("S9").Value = 1.7
$ws.Range("T9").Value = 2.02
$ws.Range("W9").Value = 2.67
$ws.Range("X9").Value = 1.42
$ws.Range("Y9").Value = 1.36
$ws.Range("Z9").Value = 2.9
$ws.Range("AA9").Value = 1.98
$ws.Range("AB9").Value = 1.75
$ws.Range("AC9").Value = 6.6
$ws.Range("AF9").Value = 8.5
$ws.Range("AG9").Value = 11.25
$ws.Range("AH9").Value = 27
$ws.Range("AI9").Value = 8
$ws.Range("AJ9").Value = 8.75
$ws.Range("AL9").Value = 90
$ws.Range("AM9").Value = 700
$ws.Range("AN9").Value = 23
$ws.Range("AO9").Value = 65
$ws.Range("AP9").Value = 26
$ws.Range("AQ9").Value = 250
$ws.Range("AR9").Value = 100
$ws.Range("AS9").Value = 80
$ws.Range("G14").Value = 5.25
$ws.Range("M14").Value = 1.03
$ws.Range("N14").Value = 10
$ws.Range("S14").Value = 1.7
$ws.Range("T14").Value = 2.1
$ws.Range("Y14").Value = 1.33
$ws.Range("Z14").Value = 3.25
$ws.Range("AA14").Value = 1.83
$ws.Range("AB14").Value = 1.83
$ws.Range("AC14").Value = 17
$ws.Range("AI14").Value = 13
$ws.Range("AM14").Value = 700
$ws.Range("AN14").Value = 7.5
$ws.Range("AO14").Value = 7.5
$ws.Range("AR14").Value = 12
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 4.33
$ws.Range("O15").Value = 1.17
$ws.Range("P15").Value = 5
$ws.Range("T15").Value = 2.4
$ws.Range("U15").Value = 1.9
$ws.Range("V15").Value = 1.95
$ws.Range("W15").Value = 2.25
$ws.Range("X15").Value = 1.57
$ws.Range("Y15").Value = 1.22
$ws.Range("Z15").Value = 3.75
$ws.Range("AA15").Value = 1.5
$ws.Range("AB15").Value = 2.5
$ws.Range("AC15").Value = 11
$ws.Range("AD15").Value = 11
$ws.Range("AI15").Value = 19
$ws.Range("AN15").Value = 19
$ws.Range("G16").Value = 3.7
$ws.Range("H16").Value = 3.6
$ws.Range("I16").Value = 2
$ws.Range("J16").Value = 4
$ws.Range("K16").Value = 2.25
$ws.Range("L16").Value = 2.6
$ws.Range("S16").Value = 1.73
$ws.Range("T16").Value = 2.08
$ws.Range("AF16").Value = 41
$ws.Range("AI16").Value = 13
$ws.Range("AJ16").Value = 7
$ws.Range("AO16").Value = 10
$ws.Range("AQ16").Value = 17
$ws.Range("AR16").Value = 15
$ws.Range("G18").Value = 1.27
$ws.Range("H18").Value = 5.1
$ws.Range("I18").Value = 9.5
$ws.Range("J18").Value = 1.72
$ws.Range("K18").Value = 2.42
$ws.Range("L18").Value = 8
$ws.Range("O18").Value = 1.2
$ws.Range("P18").Value = 3.6
$ws.Range("S18").Value = 1.62
$ws.Range("T18").Value = 2.02
$ws.Range("W18").Value = 2.47
$ws.Range("X18").Value = 1.42
$ws.Range("AA18").Value = 2.12
$ws.Range("AB18").Value = 1.57
$ws.Range("AC18").Value = 6.7
$ws.Range("AD18").Value = 5.8
$ws.Range("AE18").Value = 9
$ws.Range("AF18").Value = 7.3
$ws.Range("AG18").Value = 11.25
$ws.Range("AH18").Value = 35
$ws.Range("AI18").Value = 12.5
$ws.Range("AJ18").Value = 10.5
$ws.Range("AK18").Value = 27
$ws.Range("AL18").Value = 150
$ws.Range("AN18").Value = 23
$ws.Range("AO18").Value = 70
$ws.Range("AP18").Value = 32
$ws.Range("AQ18").Value = 300
$ws.Range("AR18").Value = 120
$ws.Range("AS18").Value = 120
